$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 01:20"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 141169
$ws.Range("C4").Value = 17591
$ws.Range("E4").Value = 134276
$ws.Range("G4").Value = 238
$ws.Range("H4").Value = 2458

# Row 8: Alemania
$ws.Range("E8").Value = 52351
$ws.Range("G8").Value = 100
$ws.Range("H8").Value = 533

# Row 17: Austria
$ws.Range("B17").Value = 8788
$ws.Range("C17").Value = 517
$ws.Range("E17").Value = 8223

# Row 20: Noruega
$ws.Range("B20").Value = 4284
$ws.Range("C20").Value = 269
$ws.Range("E20").Value = 4252

# Row 23: Australia
$ws.Range("B23").Value = 4163
$ws.Range("C23").Value = 528
$ws.Range("E23").Value = 3920
$ws.Range("F23").Value = 28
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 17

# Row 31: Ecuador
$ws.Range("B31").Value = 1924
$ws.Range("C31").Value = 101
$ws.Range("E31").Value = 1863
$ws.Range("G31").Value = 10
$ws.Range("H31").Value = 58

# Row 46: Panama
$ws.Range("B46").Value = 989
$ws.Range("C46").Value = 88
$ws.Range("E46").Value = 961
$ws.Range("G46").Value = 7
$ws.Range("H46").Value = 24

# Row 51: Argentina
$ws.Range("B51").Value = 820
$ws.Range("C51").Value = 75
$ws.Range("D51").Value = 91
$ws.Range("E51").Value = 709
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 20

# Row 104: Venezuela
$ws.Range("E104").Value = 77
$ws.Range("F104").Value = 6
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 3

# Row 107: Honduras
$ws.Range("E107").Value = 104
$ws.Range("G107").Value = 2
$ws.Range("H107").Value = 3

# Row 108: Estado de Palestina
$ws.Range("B108").Value = 109
$ws.Range("C108").Value = 5
$ws.Range("E108").Value = 90

# Row 128: Isla de Man
$ws.Range("B128").Value = 43
$ws.Range("C128").Value = 15
$ws.Range("D128").Value = 6
$ws.Range("E128").Value = 37

# Row 129: Kenia
$ws.Range("C129").Value = 10
$ws.Range("D129").Value = 0
$ws.Range("E129").Value = 42
$ws.Range("F129").Value = 0
$ws.Range("H129").Value = 0

# Row 130: Madagascar
$ws.Range("B130").Value = 42
$ws.Range("C130").Value = 4
$ws.Range("D130").Value = 1
$ws.Range("E130").Value = 40
$ws.Range("F130").Value = 2
$ws.Range("H130").Value = 1

# Row 131: Puerto Rico
$ws.Range("C131").Value = 13
$ws.Range("D131").Value = 0
$ws.Range("E131").Value = 39
$ws.Range("H131").Value = 0

# Row 132: Macao
$ws.Range("B132").Value = 39
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 1
$ws.Range("E132").Value = 36
$ws.Range("H132").Value = 2

# Row 133: Guatemala
$ws.Range("B133").Value = 37
$ws.Range("C133").Value = 3
$ws.Range("E133").Value = 27
$ws.Range("F133").Value = 0
$ws.Range("H133").Value = 0

# Row 134: Uganda
$ws.Range("B134").Value = 34
$ws.Range("C134").Value = 0
$ws.Range("D134").Value = 10
$ws.Range("E134").Value = 23
$ws.Range("F134").Value = 1
$ws.Range("H134").Value = 1

# Row 135: Barbados
$ws.Range("C135").Value = 3

# Row 136: Guam
$ws.Range("B136").Value = 33
$ws.Range("C136").Value = 7
$ws.Range("E136").Value = 33
$ws.Range("H136").Value = 0

# Row 137: Jamaica
$ws.Range("C137").Value = 0
$ws.Range("D137").Value = 0
$ws.Range("E137").Value = 31

# Row 138: Polinesia Francesa
$ws.Range("B138").Value = 32
$ws.Range("C138").Value = 2
$ws.Range("D138").Value = 2
$ws.Range("E138").Value = 29
$ws.Range("H138").Value = 1

# Row 139: Zambia
$ws.Range("B139").Value = 30
$ws.Range("C139").Value = 0
$ws.Range("E139").Value = 30

# Row 140: Guayana Francesa
$ws.Range("B140").Value = 29
$ws.Range("C140").Value = 1
$ws.Range("D140").Value = 0
$ws.Range("E140").Value = 29

# Row 144: Congo
$ws.Range("B144").Value = 21
$ws.Range("C144").Value = 5
$ws.Range("D144").Value = 1
$ws.Range("E144").Value = 20

# Row 145: Etiopia
$ws.Range("C145").Value = 15
$ws.Range("D145").Value = 0
$ws.Range("E145").Value = 19

# Row 156: Eritrea
$ws.Range("C156").Value = 0

# Row 157: Mongolia
$ws.Range("C157").Value = 6

# Row 164: Granada
$ws.Range("C164").Value = 0

# Row 165: Suazilandia
$ws.Range("C165").Value = 2

# Row 188: Fiyi
$ws.Range("D188").Value = 1
$ws.Range("E188").Value = 4

# Row 189: Nepal
$ws.Range("D189").Value = 2
$ws.Range("E189").Value = 3
